# Apply the Sablefish_Inputs.xlsx edits described in the commit:
# "breakthrough!!! removed biased estimates, and have relatively unbiased estimates now!"

$wb = $excel.ActiveWorkbook

# Controls sheet: n_sims (B2) 50 -> 100
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Range("B2").Value = 100

# Recruitment_Mortality sheet: sigma_rec (B4) 0.6 -> 0.8,
#   M (B5) 0.15 -> 0.1, mu_rec (B6) 3 -> 2.75
$wsRecruit = $wb.Worksheets.Item("Recruitment_Mortality")
$wsRecruit.Range("B4").Value = 0.8
$wsRecruit.Range("B5").Value = 0.1
$wsRecruit.Range("B6").Value = 2.75

# Move the active selection on that sheet from B5 to B7, matching the author's
# final cursor position after editing.
$wsRecruit.Activate()
$wsRecruit.Range("B7").Select()

# Restore Controls as the active/selected sheet tab (unchanged in the diff).
$wsControls.Activate()
